$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Information")

# Helper cell used to push date-looking strings ("12/12/1212",
# "12-23-2222") into cells as plain text, instead of letting Excel's
# smart typing convert them into date serial numbers. We build the
# literal text via a formula (so it is never auto-detected as a date),
# copy its value with PasteSpecial (values only), then clean up the
# helper cell so it leaves no trace in the saved workbook.
$helper = $ws.Cells.Item(1000, 1)

function Set-PlainText($cell, $text) {
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $cell.PasteSpecial(-4163)
    $helper.Clear()
}

# New row 4: a "test" account
$ws.Cells.Item(4, 1).Value = "test"
$ws.Cells.Item(4, 2).Value = "test"
$ws.Cells.Item(4, 3).Value = "test@gmail.com"
$ws.Cells.Item(4, 4).Value = "test"
Set-PlainText $ws.Cells.Item(4, 5) "12/12/1212"
$ws.Cells.Item(4, 6).Value = 69176879

# New row 5: a "notadmin" account
$ws.Cells.Item(5, 1).Value = "notadmin"
$ws.Cells.Item(5, 2).Value = "notadmin"
$ws.Cells.Item(5, 3).Value = "it workssss"
$ws.Cells.Item(5, 4).Value = "jack"
Set-PlainText $ws.Cells.Item(5, 5) "12-23-2222"
$ws.Cells.Item(5, 6).Value = 87999924
